# Apply cryptos list update per upstream data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.216.74"
$ws.Range("E2").Value = "  -0.01%  "
$ws.Range("D3").Value = "3.565.13"
$ws.Range("E3").Value = "  +0.63%  "
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E4").Value = "  -0.08%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "605.73"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.34%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "144.51"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -0.09%  "
$ws.Range("D7").Value = "3.563.66"
$ws.Range("E7").Value = "  +0.60%  "
$ws.Range("E8").Value = "  +0.19%  "
$ws.Range("E9").Value = "  +2.16%  "
$ws.Range("E10").Value = "  -0.43%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "7.81"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -2.83%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.414"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -0.27%  "
$ws.Range("D13").Value = "4.168.15"
$ws.Range("E13").Value = "  +0.63%  "
$ws.Range("E14").Value = "  -1.00%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "30.36"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -0.40%  "
$ws.Range("D16").Value = "3.554.98"
$ws.Range("E16").Value = "  +0.41%  "
$ws.Range("D17").Value = "66.265.41"
$ws.Range("E17").Value = "  -0.15%  "
$ws.Range("E18").Value = "  -0.57%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "11.42"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +4.43%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "6.23"
$c.Style = "Normal"
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "14.82"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -1.42%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "430.60"
$c.Style = "Normal"
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "0.614"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +1.79%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "79.52"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +1.22%  "
$ws.Range("D25").Value = "3.707.00"
$ws.Range("E25").Value = "  +0.71%  "
$ws.Range("E26").Value = "  +0.05%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "0.0000118"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -2.47%  "
$ws.Range("E28").Value = "  +0.73%  "
$ws.Range("E29").Value = "  -2.29%  "
$ws.Range("E30").Value = "  -1.57%  "
$ws.Range("E31").Value = "  -0.08%  "
$ws.Range("D32").Value = "3.558.77"
$ws.Range("E32").Value = "  +0.75%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "25.45"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +0.18%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "1.45"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -2.68%  "
$ws.Range("E35").Value = "  -4.33%  "
$ws.Range("B36").Value = "Aptos"
$ws.Range("C36").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "7.83"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -0.52%  "
$ws.Range("B37").Value = "USDe"
$ws.Range("C37").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +0.04%  "
$ws.Range("E38").Value = "  -1.86%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "5.61"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -0.39%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "175.47"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +3.37%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.0850"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -1.53%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "5.21"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +0.26%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.888"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -0.70%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "46.00"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +1.60%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "2.49"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +2.61%  "
$ws.Range("E48").Value = "  -2.25%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "24.94"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -4.90%  "
$ws.Range("E50").Value = "  -0.96%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "23.41"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +3.46%  "
